# Auto-applied numeric updates to Lamia_Profits market-price columns (H-N)
# per sheet: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 625882.9399999999
$ws.Range("I6").Value = 1250334.6
$ws.Range("J6").Value = 1431.25
$ws.Range("K6").Value = 3751003.8
$ws.Range("L6").Value = 4293.75
$ws.Range("M6").Value = -3750891.8
$ws.Range("N6").Value = -4517.75

$ws.Range("H19").Value = 759.6
$ws.Range("J19").Value = 849.8333
$ws.Range("L19").Value = 849.8333
$ws.Range("N19").Value = -1199.8333

$ws.Range("H62").Value = 5881.5625
$ws.Range("I62").Value = 2762.375
$ws.Range("J62").Value = 9000.75
$ws.Range("K62").Value = 2762.375
$ws.Range("L62").Value = 9000.75
$ws.Range("M62").Value = -2138.375
$ws.Range("N62").Value = -10248.75

$ws.Range("H65").Value = 5881.5625
$ws.Range("I65").Value = 2762.375
$ws.Range("J65").Value = 9000.75
$ws.Range("K65").Value = 13811.875
$ws.Range("L65").Value = 45003.75
$ws.Range("M65").Value = -10691.875
$ws.Range("N65").Value = -51243.75

$ws.Range("H100").Value = 10248.6
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H132").Value = 1643.3529
$ws.Range("I132").Value = 1541.1555
$ws.Range("K132").Value = 4623.4665
$ws.Range("M132").Value = -2093.4665

$ws.Range("H137").Value = 5214.5713
$ws.Range("I137").Value = 5241
$ws.Range("K137").Value = 15723
$ws.Range("M137").Value = -13173

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6051.079
$ws.Range("I32").Value = 4614.933
$ws.Range("K32").Value = 4614.933
$ws.Range("M32").Value = -4327.933

$ws.Range("H50").Value = 7666.6665
$ws.Range("I50").Value = 10000
$ws.Range("J50").Value = 7200
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 7200
$ws.Range("M50").Value = -9286
$ws.Range("N50").Value = -8628

$ws.Range("H109").Value = 105416.664
$ws.Range("J109").Value = 105416.664
$ws.Range("L109").Value = 105416.664
$ws.Range("N109").Value = -108190.664

$ws.Range("H131").Value = 44998.668
$ws.Range("J131").Value = 44998.668
$ws.Range("L131").Value = 44998.668
$ws.Range("N131").Value = -55078.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2631.4443
$ws.Range("I99").Value = 2637.8
$ws.Range("K99").Value = 2637.8
$ws.Range("M99").Value = -1139.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 10159.444
$ws.Range("J62").Value = 16864
$ws.Range("L62").Value = 16864
$ws.Range("N62").Value = -18112

$ws.Range("H65").Value = 10159.444
$ws.Range("J65").Value = 16864
$ws.Range("L65").Value = 84320
$ws.Range("N65").Value = -90560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1963.1333
$ws.Range("J34").Value = 2483.3333
$ws.Range("L34").Value = 7449.999899999999
$ws.Range("N34").Value = -7617.999899999999

$ws.Range("H39").Value = 2218
$ws.Range("J39").Value = 2599.3333
$ws.Range("L39").Value = 7797.999899999999
$ws.Range("N39").Value = -8385.999899999999

$ws.Range("H55").Value = 1787.375
$ws.Range("I55").Value = 1333.4166
$ws.Range("J55").Value = 3149.25
$ws.Range("K55").Value = 4000.2498
$ws.Range("L55").Value = 9447.75
$ws.Range("M55").Value = -3823.2498
$ws.Range("N55").Value = -9801.75

$ws.Range("H131").Value = 7793028
$ws.Range("I131").Value = 17857826
$ws.Range("J131").Value = 5720864
$ws.Range("K131").Value = 53573478
$ws.Range("L131").Value = 17162592
$ws.Range("M131").Value = -53568438
$ws.Range("N131").Value = -17172672

$ws.Range("H132").Value = 3916.1177
$ws.Range("J132").Value = 4320.1
$ws.Range("L132").Value = 38880.9
$ws.Range("N132").Value = -43940.9

$ws.Range("H140").Value = 2011.8823
$ws.Range("I140").Value = 1725.5
$ws.Range("K140").Value = 5176.5
$ws.Range("M140").Value = 3.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 29100.334
$ws.Range("I99").Value = 17978.166
$ws.Range("J99").Value = 51344.668
$ws.Range("K99").Value = 17978.166
$ws.Range("L99").Value = 51344.668
$ws.Range("M99").Value = -15732.166
$ws.Range("N99").Value = -55836.668

$ws.Range("H101").Value = 60000
$ws.Range("J101").Value = 60000
$ws.Range("L101").Value = 60000
$ws.Range("N101").Value = -66490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7109.0195
$ws.Range("I7").Value = 5757.054
$ws.Range("K7").Value = 5757.054
$ws.Range("M7").Value = -5645.054

$ws.Range("H22").Value = 5463.1055
$ws.Range("I22").Value = 2915.8333
$ws.Range("K22").Value = 2915.8333
$ws.Range("M22").Value = -2620.8333

$ws.Range("H27").Value = 5463.1055
$ws.Range("I27").Value = 2915.8333
$ws.Range("K27").Value = 2915.8333
$ws.Range("M27").Value = -2808.8333

$ws.Range("H46").Value = 3966.8096
$ws.Range("I46").Value = 2800.3333
$ws.Range("J46").Value = 4433.4
$ws.Range("K46").Value = 2800.3333
$ws.Range("L46").Value = 4433.4
$ws.Range("M46").Value = -2612.3333
$ws.Range("N46").Value = -4809.4

$ws.Range("H55").Value = 2502126
$ws.Range("J55").Value = 4940.857
$ws.Range("L55").Value = 4940.857
$ws.Range("N55").Value = -5286.857

$ws.Range("H61").Value = 3578.8928
$ws.Range("J61").Value = 6126.125
$ws.Range("L61").Value = 6126.125
$ws.Range("N61").Value = -6530.125

$ws.Range("H113").Value = 3578.8928
$ws.Range("J113").Value = 6126.125
$ws.Range("L113").Value = 6126.125
$ws.Range("N113").Value = -10466.125

$ws.Range("H122").Value = 171117.05
$ws.Range("I122").Value = 237875.7
$ws.Range("J122").Value = 8988.857
$ws.Range("K122").Value = 713627.1000000001
$ws.Range("L122").Value = 26966.571
$ws.Range("M122").Value = -711177.1000000001
$ws.Range("N122").Value = -31866.571

$ws.Range("H126").Value = 7109.0195
$ws.Range("I126").Value = 5757.054
$ws.Range("K126").Value = 17271.162
$ws.Range("M126").Value = -14801.162

$ws.Range("H132").Value = 5475.302
$ws.Range("I132").Value = 4854
$ws.Range("J132").Value = 6500.45
$ws.Range("K132").Value = 14562
$ws.Range("L132").Value = 19501.35
$ws.Range("M132").Value = -12032
$ws.Range("N132").Value = -24561.35

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 800
$ws.Range("I96").Value = 800
$ws.Range("K96").Value = 800
$ws.Range("M96").Value = 573

$ws.Range("H132").Value = 3170.6904
$ws.Range("J132").Value = 10981
$ws.Range("L132").Value = 32943
$ws.Range("N132").Value = -38003
